$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Ensure the updated cells keep being stored as text (matching the existing
# text-based numeric entries elsewhere in the sheet) rather than being
# auto-converted to numbers.
$ws.Range("B13:D13").NumberFormat = "@"
$ws.Range("C14:D14").NumberFormat = "@"
$ws.Range("B16:D16").NumberFormat = "@"

# Row 13: Enterprises density (per 1000 people)
$ws.Range("B13").Value = "10.06"
$ws.Range("C13").Value = "1.53"
$ws.Range("D13").Value = "11.58"

# Row 14: Employment (% of total)
$ws.Range("C14").Value = "39.46"
$ws.Range("D14").Value = "59.96"

# Row 16: Enterprises (% of total)
$ws.Range("B16").Value = "86.31"
$ws.Range("C16").Value = "13.11"
$ws.Range("D16").Value = "99.42"
